$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40/41 swap: Algorand moves to row 40, VeChain moves to row 41,
# each carrying updated D (Price) and E (Volume) values.
# D values that look numeric are prefixed with a single quote so Excel
# keeps them as text (matching the sheets existing string-typed data),
# exactly as the source data feed stores e.g. "0.2089" as text, not a number.
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2089"
$ws.Range("E40").Value = "  +1.94%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.02231"
$ws.Range("E41").Value = "  +1.76%  "

# Price (D) and Volume(1h) (E) updates for all other rows.
$updates = @(
    @{Row=2; D="25.049.89"; E="  +2.31%  "},
    @{Row=3; D="1.676.82"; E="  +1.22%  "},
    @{Row=4; D="'0.9999"; E="  -0.59%  "},
    @{Row=5; D="'331.15"; E="  +7.84%  "},
    @{Row=6; D="'0.9994"; E="  -0.09%  "},
    @{Row=7; D="'0.3662"; E="  +1.28%  "},
    @{Row=8; D="'47.43"; E="  +0.06%  "},
    @{Row=9; D="'0.3240"; E="  -0.05%  "},
    @{Row=10; D="'1.148"; E="  +2.72%  "},
    @{Row=11; D="'0.07165"; E="  +2.65%  "},
    @{Row=12; D="'1.001"; E="  +0.13%  "},
    @{Row=13; D="'6.088"; E="  +3.59%  "},
    @{Row=14; D="'19.69"; E="  +1.63%  "},
    @{Row=15; D="1.668.13"; E="  +1.01%  "},
    @{Row=16; D="'6.673"; E="  +1.98%  "},
    @{Row=17; D="'0.00001048"; E="  +0.45%  "},
    @{Row=18; D="'0.06563"; E="  +0.65%  "},
    @{Row=19; D="'0.9988"; E="  -0.14%  "},
    @{Row=20; D="'79.01"; E="  +3.71%  "},
    @{Row=21; D="'15.87"; E="  +1.52%  "},
    @{Row=22; D="'5.908"; E="  -0.07%  "},
    @{Row=23; D="'12.89"; E="  +2.38%  "},
    @{Row=24; D="24.986.14"; E="  +2.04%  "},
    @{Row=25; D="'2.449"; E="  -0.45%  "},
    @{Row=26; D="'2.382"; E="  +3.50%  "},
    @{Row=27; D="'149.02"; E="  +1.56%  "},
    @{Row=28; D="'18.76"; E="  +2.13%  "},
    @{Row=29; D="1.854.90"; E="  +1.07%  "},
    @{Row=30; D="'126.24"; E="  +2.08%  "},
    @{Row=31; D="'1.201"; E="  +0.82%  "},
    @{Row=32; D="'4.084"; E="  +2.36%  "},
    @{Row=33; D="'5.809"; E="  +2.93%  "},
    @{Row=34; D="'0.08458"; E="  +1.22%  "},
    @{Row=35; D="'1.668"; E="  -1.74%  "},
    @{Row=36; D="'12.36"; E="  +0.14%  "},
    @{Row=37; D="'5.166"; E="  -0.27%  "},
    @{Row=38; D="'0.06073"; E="  +0.53%  "},
    @{Row=39; D="'1.230"; E="  +2.40%  "},
    @{Row=42; D="'8.275"; E="  +1.08%  "},
    @{Row=43; D="'0.9988"; E="  -0.14%  "},
    @{Row=44; D="'0.5959"; E="  +0.98%  "},
    @{Row=45; D="'13.65"; E="  +8.08%  "},
    @{Row=46; D="'3.836"; E="  +2.84%  "},
    @{Row=47; D="'0.5740"; E="  +2.86%  "},
    @{Row=48; D="'124.13"; E="  +1.82%  "},
    @{Row=49; D="'1.969"; E="  +2.07%  "},
    @{Row=50; D="'0.07016"; E="  +1.67%  "},
    @{Row=51; D="'1.191"; E="  +3.23%  "}
)

foreach ($u in $updates) {
    $ws.Range("D" + $u.Row).Value = $u.D
    $ws.Range("E" + $u.Row).Value = $u.E
}
